$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (Word places this at the last edit
#    location; it currently sits inside the "Chairs at front..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Rename "Homework01.java" -> "DeepDive01.java" (the actual content edit;
#    matches commit message "Updated L1,L2,L3 with new deepdive lesson names").
$rng = $d.Content
$found = $rng.Find.Execute("Homework01.java", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "DeepDive01.java", 2)

# 3. Re-create "_GoBack" right after the new text (Word always leaves it at the
#    spot of the most recent edit). Adding a bookmark with a zero-width range
#    exactly at "end-of-paragraph-content" lands incorrectly in this engine, so
#    insert a throwaway marker character first, bookmark around it, then erase
#    the marker -- leaving a correctly collapsed bookmark in place.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("DeepDive01.java")
$rng2.Collapse(0)
$rng2.InsertAfter("~")
$d.Bookmarks.Add("_GoBack", $rng2)
$rng2.Text = ""
